$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.39168865321855
$ws.Range("C2").Value = 7.705118187527729
$ws.Range("D2").Value = 6.124506861666825
$ws.Range("F2").Value = 30.64648970081679
$ws.Range("G2").Value = 39.86261056736877
$ws.Range("H2").Value = 16.98432677521102
$ws.Range("I2").Value = 26.97121611393189
$ws.Range("K2").Value = 9.657103836434102
$ws.Range("L2").Value = 10.77732768935635
$ws.Range("B3").Value = 12.15051307766703
$ws.Range("C3").Value = 7.661987034119713
$ws.Range("D3").Value = 6.109094150988058
$ws.Range("F3").Value = 30.62306899707377
$ws.Range("G3").Value = 39.80585660675756
$ws.Range("H3").Value = 17.01940709013206
$ws.Range("I3").Value = 27.02683374229053
$ws.Range("K3").Value = 9.492121103413655
$ws.Range("L3").Value = 10.75267311188222
$ws.Range("B4").Value = 12.0028815433371
$ws.Range("C4").Value = 7.634845639825076
$ws.Range("D4").Value = 6.099386743446866
$ws.Range("F4").Value = 30.61632872983111
$ws.Range("G4").Value = 39.78225839286551
$ws.Range("H4").Value = 17.04381535474398
$ws.Range("I4").Value = 27.06616912165674
$ws.Range("K4").Value = 9.391494338371452
$ws.Range("L4").Value = 10.73972205923505
$ws.Range("B5").Value = 11.94292202951443
$ws.Range("C5").Value = 7.623622057516834
$ws.Range("D5").Value = 6.09537015496533
$ws.Range("F5").Value = 30.6155046173716
$ws.Range("G5").Value = 39.77547368209394
$ws.Range("H5").Value = 17.0544823484231
$ws.Range("I5").Value = 27.08349989781946
$ws.Range("K5").Value = 9.350713472294217
$ws.Range("L5").Value = 10.73499805262949
$ws.Range("B6").Value = 11.93298056898979
$ws.Range("C6").Value = 7.621748600788412
$ws.Range("D6").Value = 6.094699552530492
$ws.Range("F6").Value = 30.61548390761527
$ws.Range("G6").Value = 39.77451815400573
$ws.Range("H6").Value = 17.05629707707906
$ws.Range("I6").Value = 27.08645617003394
$ws.Range("K6").Value = 9.343957094858984
$ws.Range("L6").Value = 10.73424717320855
$ws.Range("B7").Value = 12.00207197503704
$ws.Range("C7").Value = 7.634694932467726
$ws.Range("D7").Value = 6.099332819484318
$ws.Range("F7").Value = 30.6163098302733
$ws.Range("G7").Value = 39.78215542422966
$ws.Range("H7").Value = 17.04395629791969
$ws.Range("I7").Value = 27.06639758553107
$ws.Range("K7").Value = 9.390943366626695
$ws.Range("L7").Value = 10.73965610336146
$ws.Range("B8").Value = 12.30849003528283
$ws.Range("C8").Value = 7.690384700814167
$ws.Range("D8").Value = 6.119243544614517
$ws.Range("F8").Value = 30.63682981315702
$ws.Range("G8").Value = 39.84070976647707
$ws.Range("H8").Value = 16.99582649292301
$ws.Range("I8").Value = 26.98931513565116
$ws.Range("K8").Value = 9.600111554185156
$ws.Range("L8").Value = 10.76837526911453
$ws.Range("B9").Value = 12.90917427315382
$ws.Range("C9").Value = 7.794263899689346
$ws.Range("D9").Value = 6.156327766211261
$ws.Range("F9").Value = 30.7375667930127
$ws.Range("G9").Value = 40.0445322822459
$ws.Range("H9").Value = 16.92424740913268
$ws.Range("I9").Value = 26.87942081755584
$ws.Range("K9").Value = 10.01319734593842
$ws.Range("L9").Value = 10.84186114509493
$ws.Range("B10").Value = 13.34548664680622
$ws.Range("C10").Value = 7.867203874272128
$ws.Range("D10").Value = 6.182346458574306
$ws.Range("F10").Value = 30.84820354569441
$ws.Range("G10").Value = 40.24801457289677
$ws.Range("H10").Value = 16.88561191927035
$ws.Range("I10").Value = 26.82398526990774
$ws.Range("K10").Value = 10.31527566298891
$ws.Range("L10").Value = 10.90603911339534
$ws.Range("B11").Value = 13.54191810794798
$ws.Range("C11").Value = 7.899624761392799
$ws.Range("D11").Value = 6.19390999739595
$ws.Range("F11").Value = 30.90640212885151
$ws.Range("G11").Value = 40.35208219319081
$ws.Range("H11").Value = 16.87107495687437
$ws.Range("I11").Value = 26.80428762844544
$ws.Range("K11").Value = 10.45174653290226
$ws.Range("L11").Value = 10.93737840456982
$ws.Range("B12").Value = 13.61593132830035
$ws.Range("C12").Value = 7.911790071391119
$ws.Range("D12").Value = 6.198249038938418
$ws.Range("F12").Value = 30.92956221959502
$ws.Range("G12").Value = 40.39312510402098
$ws.Range("H12").Value = 16.86600771507476
$ws.Range("I12").Value = 26.79762419427225
$ws.Range("K12").Value = 10.50323740369169
$ws.Range("L12").Value = 10.94954755300075
$ws.Range("B13").Value = 13.6000089330348
$ws.Range("C13").Value = 7.90917506508923
$ws.Range("D13").Value = 6.197316329742179
$ws.Range("F13").Value = 30.92452457452385
$ws.Range("G13").Value = 40.38421340820231
$ws.Range("H13").Value = 16.86707956501153
$ws.Range("I13").Value = 26.79902387183087
$ws.Range("K13").Value = 10.49215706763115
$ws.Range("L13").Value = 10.946913405355
$ws.Range("B14").Value = 13.54801508010736
$ws.Range("C14").Value = 7.900627863181973
$ws.Range("D14").Value = 6.194267774617404
$ws.Range("F14").Value = 30.90828510886775
$ws.Range("G14").Value = 40.35542617520331
$ws.Range("H14").Value = 16.87064929850746
$ws.Range("I14").Value = 26.80372346852101
$ws.Range("K14").Value = 10.45598676666929
$ws.Range("L14").Value = 10.93837356249849
$ws.Range("B15").Value = 13.51611681985542
$ws.Range("C15").Value = 7.895377824760942
$ws.Range("D15").Value = 6.192395242135045
$ws.Range("F15").Value = 30.8984837020446
$ws.Range("G15").Value = 40.33800543932561
$ws.Range("H15").Value = 16.87289286815211
$ws.Range("I15").Value = 26.8067057718386
$ws.Range("K15").Value = 10.43380547854453
$ws.Range("L15").Value = 10.93318173747378
$ws.Range("B16").Value = 13.33260152102606
$ws.Range("C16").Value = 7.865069587407334
$ws.Range("D16").Value = 6.181585216385121
$ws.Range("F16").Value = 30.84455765822473
$ws.Range("G16").Value = 40.24144338273226
$ws.Range("H16").Value = 16.88662314512538
$ws.Range("I16").Value = 26.8253838015569
$ws.Range("K16").Value = 10.30633346526024
$ws.Range("L16").Value = 10.90403360321156
$ws.Range("B17").Value = 13.21944023154547
$ws.Range("C17").Value = 7.84628004995419
$ws.Range("D17").Value = 6.174883369399231
$ws.Range("F17").Value = 30.81348463837548
$ws.Range("G17").Value = 40.1851389041414
$ws.Range("H17").Value = 16.89582496229424
$ws.Range("I17").Value = 26.83825742897938
$ws.Range("K17").Value = 10.22785372193607
$ws.Range("L17").Value = 10.88669675754346
$ws.Range("B18").Value = 13.15416427444587
$ws.Range("C18").Value = 7.835401430263027
$ws.Range("D18").Value = 6.171003033339266
$ws.Range("F18").Value = 30.79635366013417
$ws.Range("G18").Value = 40.15383805743472
$ws.Range("H18").Value = 16.90140357100391
$ws.Range("I18").Value = 26.8461815006095
$ws.Range("K18").Value = 10.18262814790491
$ws.Range("L18").Value = 10.87692733564361
$ws.Range("B19").Value = 13.13203296809232
$ws.Range("C19").Value = 7.83170593477892
$ws.Range("D19").Value = 6.169684838986068
$ws.Range("F19").Value = 30.7906810232035
$ws.Range("G19").Value = 40.14342684955908
$ws.Range("H19").Value = 16.90334148545733
$ws.Range("I19").Value = 26.84895361983588
$ws.Range("K19").Value = 10.16730236939058
$ws.Range("L19").Value = 10.87365451555964
$ws.Range("B20").Value = 13.23150652831433
$ws.Range("C20").Value = 7.84828763357465
$ws.Range("D20").Value = 6.175599447696889
$ws.Range("F20").Value = 30.81671575016576
$ws.Range("G20").Value = 40.1910205530391
$ws.Range("H20").Value = 16.89481581153823
$ws.Range("I20").Value = 26.83683322873377
$ws.Range("K20").Value = 10.23621731731049
$ws.Range("L20").Value = 10.88852140717287
$ws.Range("B21").Value = 13.56329759450513
$ws.Range("C21").Value = 7.903141436974489
$ws.Range("D21").Value = 6.195164294880725
$ws.Range("F21").Value = 30.91302468262836
$ws.Range("G21").Value = 40.36383748350762
$ws.Range("H21").Value = 16.86958890047794
$ws.Range("I21").Value = 26.80232147656514
$ws.Range("K21").Value = 10.46661635433566
$ws.Range("L21").Value = 10.94087379069969
$ws.Range("B22").Value = 13.77794676798173
$ws.Range("C22").Value = 7.93833891032853
$ws.Range("D22").Value = 6.20771866384083
$ws.Range("F22").Value = 30.98249967086585
$ws.Range("G22").Value = 40.48630114096753
$ws.Range("H22").Value = 16.85565249501563
$ws.Range("I22").Value = 26.78440427070682
$ws.Range("K22").Value = 10.61607922997272
$ws.Range("L22").Value = 10.97684386275867
$ws.Range("B23").Value = 13.66360947038547
$ws.Range("C23").Value = 7.919613894701182
$ws.Range("D23").Value = 6.201039632804503
$ws.Range("F23").Value = 30.94482563597331
$ws.Range("G23").Value = 40.42007616934029
$ws.Range("H23").Value = 16.86285702628875
$ws.Range("I23").Value = 26.7935420922667
$ws.Range("K23").Value = 10.53642665144518
$ws.Range("L23").Value = 10.95748772980718
$ws.Range("B24").Value = 13.2260520287742
$ws.Range("C24").Value = 7.847380242396802
$ws.Range("D24").Value = 6.175275793838175
$ws.Range("F24").Value = 30.81525267949245
$ws.Range("G24").Value = 40.18835812710926
$ws.Range("H24").Value = 16.89527115037528
$ws.Range("I24").Value = 26.83747548109637
$ws.Range("K24").Value = 10.23243646367707
$ws.Range("L24").Value = 10.88769586656976
$ws.Range("B25").Value = 12.74721132023897
$ws.Range("C25").Value = 7.766745931865934
$ws.Range("D25").Value = 6.146508707355546
$ws.Range("F25").Value = 30.7038584516244
$ws.Range("G25").Value = 39.97990771226733
$ws.Range("H25").Value = 16.94116444342559
$ws.Range("I25").Value = 26.90471628185504
$ws.Range("K25").Value = 9.901470083986963
$ws.Range("L25").Value = 10.82017010969083
